# Applies the "marksheet uploading" edit:
#  - Fixes a typo in Sheet1!E11 (123354 -> 12354)
#  - Adds a new student "Shami Khan" as row 12 on Sheet1
#  - Un-shares the C4:C11 "full name" formulas (each becomes its own formula)
#  - Copies all student rows (now including Shami Khan) onto Sheet2
#  - Updates the selections shown on Sheet1 and Sheet2

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Fix the roll-number typo on Sheet1 row 11 (Pulkit Aggarwal)
# ---------------------------------------------------------------------------
$ws1.Range("E11").Value = 12354

# ---------------------------------------------------------------------------
# 2. Re-enter the "full name" helper formulas on Sheet1 (rows 4-11) so that
#    each cell carries its own formula instead of one shared formula group.
# ---------------------------------------------------------------------------
for ($i = 4; $i -le 11; $i++) {
    $ws1.Range("C$i").Formula = "=A$i&"" ""&B$i"
}

# ---------------------------------------------------------------------------
# 3. Append the new student record (Shami Khan) as row 12 on Sheet1
# ---------------------------------------------------------------------------
$row12 = New-Object 'object[,]' 1,12
$row12[0,0]  = "Shami"
$row12[0,1]  = "Khan"
$row12[0,2]  = $null
$row12[0,3]  = "X"
$row12[0,4]  = 12355
$row12[0,5]  = "Amjad Khan"
$row12[0,6]  = 9874123654
$row12[0,7]  = "Ameena Khan"
$row12[0,8]  = "855 Kachla Baba"
$row12[0,9]  = "Lamdapur"
$row12[0,10] = "Delhi"
$row12[0,11] = 110001
$ws1.Range("A12:L12").Value = $row12
$ws1.Range("C12").Formula = "=A12&"" ""&B12"

# ---------------------------------------------------------------------------
# 4. Update Sheet1's shown selection to A2:L12
# ---------------------------------------------------------------------------
$ws1.Activate()
[void]$ws1.Range("A2:L12").Select()

# ---------------------------------------------------------------------------
# 5. Populate Sheet2 with every student record (rows 1-11), mirroring
#    Sheet1 rows 2-12, including Shami Khan.
# ---------------------------------------------------------------------------
$data = @(
    ,("Sultan",   "Ahmad",    "X", 12345, "Raza Ahmad",           9654123123, "Halima Raza",          "2681 Chacha walan",  "Near Gadha Park", "Delhi", 110018)
    ,("Kashif",   "Khan",     "X", 12346, "Zalim Khan",           9652123123, "Asfa Zalim",           "142 Hata Kala",      "Lal Kuan",         "Delhi", 110006)
    ,("Rubeen",   "Hatoon",   "X", 12347, "Kashif Sultan",        9212392123, "Asma Hatun",           "125 Gali Pahar",     "Neela Pani",       "Delhi", 110041)
    ,("Rahul",    "Sharma",   "X", 12348, "Gopesh Sharma",        9658789789, "Seema Sharma",         "856 Kalua street",   "Patak Tel",        "Delhi", 110022)
    ,("Kajal",    "Gupta",    "X", 12349, "Rajat Gupta",          9632174589, "Anju Gupta",           "947 Jamshed Nagar",  "Lala Pare",        "Delhi", 110011)
    ,("Zubaida",  "Khan",     "X", 12350, "Salman Paasha",        9212654654, "Latifa Hatun",         "966 Ahmed Nagar",    "Lal mandi",        "Delhi", 110033)
    ,("Rajat",    "Sharma",   "X", 12351, "Gaurav Sharma",        9632587417, "Sameera Sharma",       "66 Gali Jor",        "Hata baazar",      "Delhi", 110045)
    ,("Ankit",    "Mishra",   "X", 12352, "BK Mishra",            9685236514, "Anita Mishra",         "1 Kali mata road",   "janda bazar",      "Delhi", 110025)
    ,("Zoya",     "Ahmad",    "X", 12353, "Hammad Ahmad",         9652123123, "Shabana Khatoon",      "145 Chatta Sheikh",  "Jang Pura",        "Delhi", 110042)
    ,("Pulkit",   "Aggarwal", "X", 12354, "Rahul Aggarwal",       9213456859, "Rajni Rahul Aggarwal", "989 Bazar Sita Ram", "Chawri Bazar",     "Delhi", 110006)
    ,("Shami",    "Khan",     "X", 12355, "Amjad Khan",           9874123654, "Ameena Khan",          "855 Kachla Baba",    "Lamdapur",         "Delhi", 110001)
)

$arr = New-Object 'object[,]' 11,12
for ($r = 0; $r -lt 11; $r++) {
    $rec = $data[$r]
    $arr[$r,0]  = $rec[0]
    $arr[$r,1]  = $rec[1]
    $arr[$r,2]  = $null
    $arr[$r,3]  = $rec[2]
    $arr[$r,4]  = $rec[3]
    $arr[$r,5]  = $rec[4]
    $arr[$r,6]  = $rec[5]
    $arr[$r,7]  = $rec[6]
    $arr[$r,8]  = $rec[7]
    $arr[$r,9]  = $rec[8]
    $arr[$r,10] = $rec[9]
    $arr[$r,11] = $rec[10]
}
$ws2.Range("A1:L11").Value = $arr

for ($i = 1; $i -le 11; $i++) {
    $ws2.Range("C$i").Formula = "=A$i&"" ""&B$i"
}

# ---------------------------------------------------------------------------
# 6. Update Sheet2's shown selection to A1:L11
# ---------------------------------------------------------------------------
$ws2.Activate()
[void]$ws2.Range("A1:L11").Select()

$ws1.Activate()
